$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: insert a new "Meta description" paragraph right after the H1 title
#         paragraph ("Play Blazing X free: Review of unique oriental-themed
#         slot"), before the "GAMEPLAY MECHANICS" Heading2 paragraph.
#
# We build the insertion as two <w:p> elements: the new paragraph itself,
# followed by a verbatim copy of the existing "GAMEPLAY MECHANICS" paragraph.
# InsertXML() replaces the target range's contents, so targeting it at the
# whole "GAMEPLAY MECHANICS" paragraph range and re-supplying that paragraph
# unchanged as the trailing fragment effectively prepends our new paragraph
# in front of it without disturbing anything else.
# ---------------------------------------------------------------------------

$headingPara = $d.Paragraphs.Item(2)
$headingRange = $headingPara.Range

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$metaParaXml = '<w:p ' + $wNs + '>' +
  '<w:r/>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
  '<w:r><w:t>: Experience gameplay mechanics with Wild center reel and free spins with unlimited multiplier in Blazing X. Play free and read our review.</w:t></w:r>' +
  '</w:p>'

$headingParaXml = '<w:p ' + $wNs + '>' +
  '<w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' +
  '<w:r><w:t>GAMEPLAY MECHANICS</w:t></w:r>' +
  '</w:p>'

$headingRange.InsertXML($metaParaXml + $headingParaXml)

# ---------------------------------------------------------------------------
# Part 2: near the end of the document, drop the bold duplicate-title
#         paragraph ("Play Blazing X free: Review of unique oriental-themed
#         slot") and rewrite the text of the following italic paragraph
#         (formerly the meta-description sentence) into the new AI image
#         prompt, keeping its run formatting (the <w:i/> run) and its
#         leading empty run intact.
# ---------------------------------------------------------------------------

$paraCount = $d.Paragraphs.Count
$boldTitlePara = $d.Paragraphs.Item($paraCount - 1)
$boldTitlePara.Range.Delete()

$paraCount = $d.Paragraphs.Count
$imagePromptPara = $d.Paragraphs.Item($paraCount)
$fullRange = $imagePromptPara.Range

# Exclude the trailing paragraph-mark character so only the run text changes.
$textOnlyRange = $d.Range($fullRange.Start, $fullRange.End - 1)
$textOnlyRange.Text = 'Create a feature image for "Blazing X". The image should be in a cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be holding a dragon and standing in front of the Blazing X game grid. The background should be highlighted in fiery orange and red hues to represent the theme of the game. Use creative illustrations and vivid colors to make the image pop and capture the attention of slot players looking for a new and exciting game to play.'

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
